# More Feature of Aspose.Slides v 1.2
# - Duplicate slide 2 ("Simple Presentation - Slide 2") to create a new
#   slide 3, then update its number text to "3".
# - Tidy up slide 2's title textbox: merge the "Presentation - " / "Slide "
#   runs into a single "Presentation - Slide " run (matches how PowerPoint
#   collapses runs with identical formatting after a retype).

$p = $ppt.ActivePresentation

# --- Slide 2: merge "Presentation - " + "Slide " runs into one run -------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(1)
$tr2 = $shp2.TextFrame.TextRange
$merged2 = $tr2.Characters(8, 21)
$merged2.Text = "Presentation " + [char]0x2013 + " Slide "
# Re-assert the textbox's original autofit height (retyping nudges the
# engine's relayout by a hair; PowerPoint itself settles back on this).
$shp2.Height = 41.25

# --- Duplicate slide 2 to create slide 3, then retarget its number -------
$newRange = $s2.Duplicate()
$s3 = $newRange.Item(1)
$shp3 = $s3.Shapes.Item(1)
$tr3 = $shp3.TextFrame.TextRange
$num3 = $tr3.Characters($tr3.Length, 1)
$num3.Text = "3"
$shp3.Height = 41.25
